$d = $word.ActiveDocument

# 1. Remove the "I recommend ImageJ over Fiji ... DiameterJ1.xxx.zip file." sentences
#    (leftover recommendation text consolidated/removed per commit message)
$d.Content.Find.Execute(
    "I recommend ImageJ over Fiji if you have no experience with either software because it is simpler to use and the plugin is buggy in Fiji.  Soon a version will be released that will work for both but for now I recommend ImageJ.  If you install imageJ version 1.49n or newer (including and especially ImageJ 2.xxx) please download the DimaeterJ2.xxx.zip file instead of DiameterJ1.xxx.zip file.  ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "",
    2
)

# 2. Note the FIJI version only needs one folder
$d.Content.Find.Execute(
    "move or copy the three folders into the plugins folder of ImageJ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "move or copy the three folders (or one folder in the FIJI version) into the plugins folder of ImageJ",
    2
)

# 3. Rename the "Segmented Images" output folder to "Segmented Images_XXX"
$d.Content.Find.Execute(
    "The “Segmented Images" + [char]34 + " folder has all of the segmented images",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The “Segmented Images_XXX" + [char]34 + " folder has all of the segmented images",
    2
)

# Find/Replace auto-corrects a straight quote that follows a letter into a curly
# closing quote ("smart quotes"); put the straight quote back to match the
# original author's (inconsistent) punctuation style.
$fixRange = $d.Content
$foundCurly = $fixRange.Find.Execute("Images_XXX" + [char]8221)
if ($foundCurly) {
    $fixRange.Text = "Images_XXX" + [char]34
}

# 4. Replace the "They are pretty self-explanatory." sentence with a pointer to the
#    new "DiameterJ Output Descriptions.docx" reference file.
$d.Content.Find.Execute(
    "They are pretty self-explanatory. ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "The file labeled “DiameterJ Output Descriptions.docx” describes what the outputs are in each of these folders.",
    2
)
